$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) cells we touch stay plain text (avoids
# Excel auto-converting strings like "227.94" into numeric values and
# dropping significant trailing zeros, e.g. "5.90" -> 5.9).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.796.32"
$ws.Range("E2").Value = "  -0.62%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.035.68"
$ws.Range("E3").Value = "  -0.65%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.94"
$ws.Range("E5").Value = "  +0.07%  "

$ws.Range("E6").Value = "  -1.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.21"
$ws.Range("E7").Value = "  -1.29%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.374"
$ws.Range("E9").Value = "  -2.32%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0827"
$ws.Range("E10").Value = "  +1.27%  "

$ws.Range("E11").Value = "  +0.16%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.336.81"
$ws.Range("E12").Value = "  -0.83%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.51"
$ws.Range("E13").Value = "  -1.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.97"
$ws.Range("E14").Value = "  -0.21%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.770"
$ws.Range("E15").Value = "  +1.58%  "

$ws.Range("E16").Value = "  -0.28%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.056.56"
$ws.Range("E17").Value = "  +0.13%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.732.27"
$ws.Range("E18").Value = "  -0.72%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.37"
$ws.Range("E19").Value = "  -0.49%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.90"
$ws.Range("E20").Value = "  -3.54%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0820"
$ws.Range("E21").Value = "  -0.76%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.75"
$ws.Range("E22").Value = "  -0.32%  "

$ws.Range("E23").Value = "  +0.14%  "

$ws.Range("E24").Value = "  -0.25%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.28"
$ws.Range("E25").Value = "  +2.98%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.10"
$ws.Range("E26").Value = "  +1.17%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.31"
$ws.Range("E27").Value = "  +1.41%  "

$ws.Range("E28").Value = "  -1.62%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.72"
$ws.Range("E29").Value = "  -1.02%  "

$ws.Range("E30").Value = "  -2.11%  "

$ws.Range("E31").Value = "  -0.10%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.20"
$ws.Range("E32").Value = "  +7.42%  "

$ws.Range("E33").Value = "  -2.10%  "

$ws.Range("E34").Value = "  -0.15%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.47"
$ws.Range("E35").Value = "  -1.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.52"
$ws.Range("E36").Value = "  +3.55%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.33"
$ws.Range("E37").Value = "  +1.51%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.39"
$ws.Range("E38").Value = "  +4.89%  "

$ws.Range("E39").Value = "  -0.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.89"
$ws.Range("E40").Value = "  +7.77%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.526.99"
$ws.Range("E41").Value = "  -0.33%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.07"
$ws.Range("E42").Value = "  -0.39%  "

$ws.Range("E43").Value = "  -0.99%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.83"
$ws.Range("E44").Value = "  -0.19%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0907"
$ws.Range("E45").Value = "  -1.96%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.14"
$ws.Range("E46").Value = "  +3.57%  "

$ws.Range("E47").Value = "  -1.06%  "

$ws.Range("E48").Value = "  -0.39%  "

$ws.Range("E49").Value = "  -1.79%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.07"
$ws.Range("E50").Value = "  +0.41%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.224.92"
$ws.Range("E51").Value = "  -0.88%  "
